$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "date"
$ws.Range("B1").Value = "week"
$ws.Range("C1").Value = "user"
$ws.Range("D1").Value = "task"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# --- Row 2 ---
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-07-21"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "週一"
$ws.Range("C2").Value = "小A"
$ws.Range("D2").Value = "教學"

# --- Row 3 ---
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2025-07-21"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "週一"
$ws.Range("C3").Value = "小B"
$ws.Range("D3").Value = "備課"

# --- Row 4 (new row) ---
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2025-07-22"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = "週二"
$ws.Range("C4").Value = "小C"
$ws.Range("D4").Value = "教學"
